# Auto-generated Excel COM-interop script
# Applies cryptos price/volume updates for Sat Mar 11 07:40:19 UTC 2023
# (also swaps two pairs of rows whose Coin/Link/Price/Volume shifted order)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings
# (e.g. "39.50", "0.00001028") keep their exact original formatting
# instead of being parsed into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '20.369.57'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '1.467.67'
$ws.Range("E3").Value = '  +3.95%  '
$ws.Range("D4").Value = '1.015'
$ws.Range("E4").Value = '  +1.63%  '
$ws.Range("D5").Value = '279.88'
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").Value = '0.8991'
$ws.Range("E6").Value = '  -10.12%  '
$ws.Range("D7").Value = '0.3717'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.3173'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").Value = '39.50'
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").Value = '1.042'
$ws.Range("E10").Value = '  +4.18%  '
$ws.Range("D11").Value = '0.06578'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.007'
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '5.508'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '17.73'
$ws.Range("E14").Value = '  +4.21%  '
$ws.Range("D15").Value = '6.169'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '1.474.32'
$ws.Range("E16").Value = '  +4.86%  '
$ws.Range("D17").Value = '0.00001028'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").Value = '0.05654'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("D19").Value = '0.9035'
$ws.Range("E19").Value = '  -9.64%  '
$ws.Range("D20").Value = '69.42'
$ws.Range("E20").Value = '  -5.98%  '
$ws.Range("D21").Value = '5.650'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '14.53'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '11.12'
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("D24").Value = '2.284'
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("D25").Value = '20.542.57'
$ws.Range("E25").Value = '  +2.68%  '
$ws.Range("D26").Value = '2.243'
$ws.Range("E26").Value = '  -1.66%  '
$ws.Range("D27").Value = '136.93'
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = '17.31'
$ws.Range("E28").Value = '  +2.35%  '
$ws.Range("D29").Value = '1.631.82'
$ws.Range("E29").Value = '  +4.29%  '
$ws.Range("D30").Value = '112.51'
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("E31").Value = '  +3.16%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.8216'
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.966'
$ws.Range("E33").Value = '  -8.04%  '
$ws.Range("D34").Value = '0.07791'
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").Value = '0.06029'
$ws.Range("E35").Value = '  +3.78%  '
$ws.Range("D36").Value = '1.472'
$ws.Range("E36").Value = '  +15.57%  '
$ws.Range("D37").Value = '4.814'
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("D38").Value = '1.157'
$ws.Range("E38").Value = '  +8.74%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02029'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '10.38'
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("D41").Value = '0.1853'
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("D42").Value = '0.9194'
$ws.Range("E42").Value = '  -8.07%  '
$ws.Range("D43").Value = '0.5317'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '3.567'
$ws.Range("E44").Value = '  +0.96%  '
$ws.Range("D45").Value = '6.884'
$ws.Range("E45").Value = '  -18.69%  '
$ws.Range("D46").Value = '12.20'
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").Value = '121.76'
$ws.Range("E47").Value = '  +10.52%  '
$ws.Range("D48").Value = '0.5226'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("D49").Value = '1.811'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").Value = '0.06407'
$ws.Range("E50").Value = '  +3.81%  '
$ws.Range("D51").Value = '1.027'
$ws.Range("E51").Value = '  -1.91%  '
